$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.033558130264282
$ws.Range("B1").Value = 1.638951659202576
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 0.3626376390457153
